$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 12

$ws.Cells.Item($row, 1).Value = 42620.891168981485
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item($row, 2).Value = 2
$ws.Cells.Item($row, 3).Value = 55
$ws.Cells.Item($row, 4).Value = 37
$ws.Cells.Item($row, 5).Value = 55
$ws.Cells.Item($row, 6).Value = 50
$ws.Cells.Item($row, 7).Value = 31869
$ws.Cells.Item($row, 8).Value = 19383
$ws.Cells.Item($row, 9).Value = 3143
$ws.Cells.Item($row, 10).Value = 461
$ws.Cells.Item($row, 11).Value = 310
$ws.Cells.Item($row, 12).Value = 2
$ws.Cells.Item($row, 13).Value = 2
$ws.Cells.Item($row, 14).Value = "Noun"
